{"js": "// Added Color Coded Sheet of Specs\n//\n// The document body is one giant paragraph that uses manual line breaks\n// (<w:br/>, represented as Chr(11)/vertical-tab in Range text) instead of\n// real paragraph marks. After each \"Utilities Req'd:\" line there are two\n// consecutive manual breaks (a blank line before \"See plans for location...\").\n// This edit turns the FIRST of those two manual breaks into a real paragraph\n// break, splitting the big paragraph in two at that point, while leaving the\n// SECOND manual break in place (so the blank-line spacing is preserved).\n//\n// We do this for both equipment blocks in the document (the \"30.5A\" block\n// and the \"31.0A\" block): locate the text immediately before the break\n// pair, insert a real paragraph mark right after the first break, then\n// delete that now-redundant first break character.\n\nconst VT = String.fromCharCode(11); // manual line break (<w:br/>) as text\nconst CR = String.fromCharCode(13); // paragraph mark when inserted as text\n\nasync function splitAfterFirstBreak(anchorText) {\n  // Find \"<anchorText><break>\" -- this anchors us right after the first of\n  // the two manual breaks that follow the Utilities line.\n  const results = context.document.body.search(anchorText + VT, { matchCase: true });\n  results.load(\"items,text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find anchor text: \" + anchorText);\n  }\n  const hit = results.items[0];\n\n  // Insert a real paragraph break right after that first manual break --\n  // this lands exactly between the two manual breaks, splitting the\n  // paragraph in two.\n  hit.insertText(CR, \"End\");\n  await context.sync();\n\n  // Re-find the same text (now at the tail of the first half paragraph) so\n  // we can clean up the now-redundant first manual break.\n  const results2 = context.document.body.search(anchorText + VT, { matchCase: true });\n  results2.load(\"items,text\");\n  await context.sync();\n  const hit2 = results2.items[0];\n\n  // Isolate just the trailing manual-break character within that hit, and\n  // delete it (the new paragraph mark now does its job).\n  const subResults = hit2.search(VT, { matchCase: true });\n  subResults.load(\"items,text\");\n  await context.sync();\n  subResults.items[0].delete();\n  await context.sync();\n}\n\nawait splitAfterFirstBreak(\"208V/3PH; (2) 30.5A\");\nawait splitAfterFirstBreak(\"(2) 208V/3PH; 31.0A\");\n", "ps1": "# Added Color Coded Sheet of Specs\n#\n# The document body is one giant paragraph that uses manual line breaks\n# (<w:br/>, rendered as Chr(11)/vertical-tab in Range.Text) instead of real\n# paragraph marks. After each \"Utilities Req'd:\" line there are two\n# consecutive manual breaks (a blank line before \"See plans for location...\").\n# This edit turns the FIRST of those two manual breaks into a real paragraph\n# break, splitting the big paragraph in two at that point, while leaving the\n# SECOND manual break in place (so the blank-line spacing is preserved).\n#\n# We do this for both equipment blocks in the document (the \"30.5A\" block and\n# the \"31.0A\" block), by locating the text immediately preceding the break\n# pair, inserting a real paragraph mark right after the first break, and then\n# deleting that first break character.\n\n$d = $word.ActiveDocument\n\nfunction Split-ParagraphAfterUtilities([string]$anchorText) {\n    # Find the line-specific text that immediately precedes the pair of\n    # manual line breaks we need to split on.\n    $rng = $d.Content\n    $found = $rng.Find.Execute($anchorText)\n    if (-not $found) {\n        throw \"Could not find anchor text '$anchorText'\"\n    }\n\n    # Collapse to the end of the match -- this is the character position\n    # right before the first manual break (Chr(11)).\n    $rng.Collapse(0)   # wdCollapseEnd\n    $breakStart = $rng.Start\n\n    # Sanity check: the next two characters should be the manual breaks.\n    $twoChars = $d.Range($breakStart, $breakStart + 2).Text\n    if ($twoChars.Length -ne 2 -or [int][char]$twoChars[0] -ne 11 -or [int][char]$twoChars[1] -ne 11) {\n        throw \"Unexpected content after '$anchorText' (expected two manual breaks)\"\n    }\n\n    # Insert a real paragraph mark right after the first manual break (i.e.\n    # right before the second one), splitting the paragraph there.\n    $insertPoint = $d.Range($breakStart + 1, $breakStart + 1)\n    $insertPoint.InsertParagraphBefore()\n\n    # Remove the first manual break run -- its job is now done by the new\n    # paragraph mark.\n    $d.Range($breakStart, $breakStart + 1).Delete()\n}\n\nSplit-ParagraphAfterUtilities \"208V/3PH; (2) 30.5A\"\nSplit-ParagraphAfterUtilities \"(2) 208V/3PH; 31.0A\"\n"}
